$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.169.19"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "3.272.64"
$ws.Range("E3").Value = "  +0.06%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "184.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "576.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.129"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.406"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").Value = "3.861.69"
$ws.Range("E12").Value = "  +0.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.137"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.76%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.24"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "67.490.16"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "3.318.31"
$ws.Range("E17").Value = "  +1.29%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "436.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +8.54%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.11"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "

$ws.Range("B24").Value = "WrappedeETH"
$ws.Range("C24").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D24").Value = "3.449.11"
$ws.Range("E24").Value = "  +1.09%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.512"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000117"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.188"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.15%  "

$ws.Range("E29").Value = "  -0.47%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.70"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.76%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.23"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.15%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.75"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.27%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.42%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.47%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.784"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.60%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.708.64"
$ws.Range("E42").Value = "  +1.22%  "

$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.35"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0670"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.03%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "24.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "324.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0272"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.42%  "

$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.988"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.18%  "

$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.27%  "
